$d = $word.ActiveDocument

# The document contains three runs whose <w:t> text ends (or is split) with a
# literal newline character that should be collapsed:
#   1. "This is an annotatable resource in the casebook.\n" -> trailing "\n" removed
#   2. ";\nreplaced: "                                       -> "\n" becomes " "
#   3. "; noted:\n"                                          -> "\n" becomes " "
#
# Using Find/Replace across the whole run merges the edited text into the
# neighboring (differently-styled) run, so instead we locate each bare
# newline character individually and fix only that single character, which
# keeps every run's formatting/boundaries intact.

function Fix-NextNewline {
    param([string]$Replacement)

    $probe = $d.Range(0, $d.Content.End)
    $found = $probe.Find.Execute("`n", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return $false
    }

    $target = $d.Range($probe.Start, $probe.End)
    $target.Text = $Replacement
    return $true
}

# 1) Drop the trailing newline after "...casebook."
Fix-NextNewline "" | Out-Null

# 2) Turn ";\nreplaced: " into "; replaced: "
Fix-NextNewline " " | Out-Null

# 3) Turn "; noted:\n" into "; noted: "
Fix-NextNewline " " | Out-Null
